$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.230.06'
$ws.Range('E2').Value = '  -2.04%  '
$ws.Range('D3').Value = '3.387.93'
$ws.Range('E3').Value = '  -1.66%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'570.37"
$ws.Range('E5').Value = '  -1.66%  '
$ws.Range('D6').Value = "'140.43"
$ws.Range('E6').Value = '  -5.91%  '
$ws.Range('D8').Value = '3.386.86'
$ws.Range('E8').Value = '  -1.71%  '
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('E10').Value = '  -4.34%  '
$ws.Range('E11').Value = '  -1.23%  '
$ws.Range('D12').Value = "'0.390"
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('D13').Value = '3.968.39'
$ws.Range('E13').Value = '  -1.62%  '
$ws.Range('D14').Value = "'28.14"
$ws.Range('E14').Value = '  +0.51%  '
$ws.Range('E15').Value = '  +0.88%  '
$ws.Range('D16').Value = '3.389.29'
$ws.Range('E16').Value = '  -1.51%  '
$ws.Range('E17').Value = '  -2.63%  '
$ws.Range('D18').Value = '60.392.65'
$ws.Range('E18').Value = '  -1.95%  '
$ws.Range('E19').Value = '  -1.39%  '
$ws.Range('D20').Value = "'14.03"
$ws.Range('E20').Value = '  -1.96%  '
$ws.Range('D21').Value = "'9.10"
$ws.Range('E21').Value = '  -4.17%  '
$ws.Range('D22').Value = "'388.55"
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').Value = "'0.560"
$ws.Range('E23').Value = '  -1.70%  '
$ws.Range('D24').Value = "'73.49"
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('E26').Value = '  -4.80%  '
$ws.Range('D27').Value = '3.527.59'
$ws.Range('E27').Value = '  -1.70%  '
$ws.Range('E28').Value = '  -1.04%  '
$ws.Range('D29').Value = "'0.998"
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  -5.43%  '
$ws.Range('D31').Value = "'7.99"
$ws.Range('E31').Value = '  -3.35%  '
$ws.Range('E32').Value = '  -1.76%  '
$ws.Range('D33').Value = "'1.42"
$ws.Range('E33').Value = '  -7.16%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').Value = "'23.72"
$ws.Range('E35').Value = '  -1.18%  '
$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.414.88'
$ws.Range('E36').Value = '  -1.57%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').Value = "'6.92"
$ws.Range('E37').Value = '  -2.29%  '
$ws.Range('D38').Value = "'168.08"
$ws.Range('E38').Value = '  +1.34%  '
$ws.Range('E39').Value = '  -6.82%  '
$ws.Range('E40').Value = '  -4.45%  '
$ws.Range('D41').Value = "'0.0775"
$ws.Range('E41').Value = '  -2.38%  '
$ws.Range('D42').Value = "'27.03"
$ws.Range('E42').Value = '  +2.01%  '
$ws.Range('E43').Value = '  -1.34%  '
$ws.Range('D44').Value = "'0.999"
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('E45').Value = '  -1.41%  '
$ws.Range('E46').Value = '  -1.48%  '
$ws.Range('D47').Value = "'41.30"
$ws.Range('E47').Value = '  -2.29%  '
$ws.Range('D48').Value = '2.520.40'
$ws.Range('E48').Value = '  -3.54%  '
$ws.Range('D49').Value = "'1.11"
$ws.Range('E49').Value = '  -3.92%  '
$ws.Range('D50').Value = "'23.14"
$ws.Range('E50').Value = '  -0.20%  '
$ws.Range('E51').Value = '  -3.52%  '
